# Updating perf results, Removing unwanted files, Refactoring
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Numbers")

# --- Row 36 (raw sample inputs for the "0.5" percentile bucket) ---
$ws.Range("L36").Value = 209
$ws.Range("M36").Value = 218
$ws.Range("N36").Value = 219
$ws.Range("O36").Value = 210

# --- Row 37 ---
$ws.Range("B37").Value = 7.7141330000000004
$ws.Range("C37").Value = 7.6513229999999997
$ws.Range("D37").Value = 7.9017299999999997
$ws.Range("E37").Value = 7.6940770000000001
$ws.Range("L37").Value = 215
$ws.Range("M37").Value = 230
$ws.Range("N37").Value = 228
$ws.Range("O37").Value = 220

# --- Row 38 ---
$ws.Range("L38").Value = 219
$ws.Range("M38").Value = 234
$ws.Range("N38").Value = 234
$ws.Range("O38").Value = 226

# --- Row 39 ---
$ws.Range("B39").Value = 27074
$ws.Range("C39").Value = 2977
$ws.Range("D39").Value = 3004
$ws.Range("E39").Value = 2968
$ws.Range("L39").Value = 221
$ws.Range("M39").Value = 237
$ws.Range("N39").Value = 238
$ws.Range("O39").Value = 229

# --- Row 40 ---
$ws.Range("L40").Value = 227
$ws.Range("M40").Value = 248
$ws.Range("N40").Value = 245
$ws.Range("O40").Value = 239

# --- Row 41 ---
$ws.Range("L41").Value = 235
$ws.Range("M41").Value = 256
$ws.Range("N41").Value = 252
$ws.Range("O41").Value = 252

# --- Row 42 ---
$ws.Range("B42").Value = 4362864
$ws.Range("C42").Value = 4362138
$ws.Range("D42").Value = 4362268
$ws.Range("E42").Value = 4362192
$ws.Range("L42").Value = 250
$ws.Range("M42").Value = 279
$ws.Range("N42").Value = 276
$ws.Range("O42").Value = 287

# --- Row 43 ---
$ws.Range("B43").Value = 1242552
$ws.Range("C43").Value = 1242138
$ws.Range("D43").Value = 1242060
$ws.Range("E43").Value = 1242192
$ws.Range("L43").Value = 387
$ws.Range("M43").Value = 1647
$ws.Range("N43").Value = 1618
$ws.Range("O43").Value = 807

# --- Row 44 ---
$ws.Range("B44").Value = 3888.97
$ws.Range("C44").Value = 3920.89
$ws.Range("D44").Value = 3796.64
$ws.Range("E44").Value = 3899.1
$ws.Range("L44").Value = 3401
$ws.Range("M44").Value = 3459
$ws.Range("N44").Value = 3444
$ws.Range("O44").Value = 3422

# --- Row 45 ---
$ws.Range("B45").Value = 257.13799999999998
$ws.Range("C45").Value = 255.04400000000001
$ws.Range("D45").Value = 263.39100000000002
$ws.Range("E45").Value = 256.46899999999999

# --- Row 46 ---
$ws.Range("B46").Value = 0.25700000000000001
$ws.Range("C46").Value = 0.255
$ws.Range("D46").Value = 0.26300000000000001
$ws.Range("E46").Value = 0.25600000000000001

# --- Row 47 ---
$ws.Range("B47").Value = 552.23
$ws.Range("C47").Value = 556.64
$ws.Range("D47").Value = 539.12
$ws.Range("E47").Value = 553.54

# Force recalculation so dependent AVERAGE/MEDIAN formulas refresh
$excel.CalculateFullRebuild()

# Restore the scrolled/selected view state on the Numbers sheet
$ws.Range("K52").Select()
